## affichage coups grille pas fini
## Adds 5 new journal rows (62-66) documenting work on displaying shots on the
## battleship grid, tweaks the H49 cell into a hyperlink (re-confirmed URL),
## updates the H55 note and expands the Tableau4 table / conditional format
## ranges to keep up with the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Grow the table (ListObject) so the new rows become part of Tableau4.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:H66")) | Out-Null

# ---------------------------------------------------------------------------
# 2. Clone row formatting (borders/number formats/fonts) onto the new rows
#    before writing any values, matching the existing alternating pattern.
# ---------------------------------------------------------------------------
$ws.Range("A61:H61").Copy() | Out-Null
$ws.Range("A62:H62").PasteSpecial(-4122) | Out-Null
$ws.Range("A55:H55").Copy() | Out-Null
$ws.Range("A63:H63").PasteSpecial(-4122) | Out-Null
$ws.Range("A61:H61").Copy() | Out-Null
$ws.Range("A64:H65").PasteSpecial(-4122) | Out-Null
$ws.Range("A55:H55").Copy() | Out-Null
$ws.Range("A66:H66").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Rows 55 and 62 pick up an explicit (but default-valued) row height.
$ws.Rows("55").RowHeight = 15
$ws.Rows("62").RowHeight = 15

# ---------------------------------------------------------------------------
# 3. Row 62 - started asking for the column letter -> number conversion.
# ---------------------------------------------------------------------------
$ws.Range("A62").Value = 44623
$ws.Range("B62").Value = 0.34375
$ws.Range("C62").Value = 0.36041666666666666
$ws.Range("D62").Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$ws.Range("E62").Value = "CPNV"
$ws.Range("F62").Value = "Demander les coordonées"
$ws.Range("G62").Value = "J'ai transformer la lettre de la colonne en numéro(ex.: A = 1)"
$ws.Range("H62").Value = "https://stackoverflow.com/questions/1469711/converting-letters-to-numbers-in-c"

# ---------------------------------------------------------------------------
# 4. Row 63 - commented the ColumnNumber function.
# ---------------------------------------------------------------------------
$ws.Range("A63").Value = 44623
$ws.Range("B63").Value = 0.36249999999999999
$ws.Range("C63").Value = 0.3659722222222222
$ws.Range("D63").Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$ws.Range("E63").Value = "CPNV"
$ws.Range("F63").Value = "Commenter"
$ws.Range("G63").Value = "J'ai commenter ma fonction ColumnNumber"

# ---------------------------------------------------------------------------
# 5. Row 64 - started displaying the shots on the grid.
# ---------------------------------------------------------------------------
$ws.Range("A64").Value = 44623
$ws.Range("B64").Value = 0.3666666666666667
$ws.Range("C64").Value = 0.39930555555555558
$ws.Range("D64").Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$ws.Range("E64").Value = "CPNV"
$ws.Range("F64").Value = "Affichage des coups dans la grille"
$ws.Range("G64").Value = "J'ai comencé la fonction qui affiche les coups dans la grille"

# ---------------------------------------------------------------------------
# 6. Row 65 - finished the function but it doesn't work yet.
# ---------------------------------------------------------------------------
$ws.Range("A65").Value = 44623
$ws.Range("B65").Value = 0.41875000000000001
$ws.Range("C65").Value = 0.44236111111111115
$ws.Range("D65").Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$ws.Range("E65").Value = "CPNV"
$ws.Range("F65").Value = "Affichage des coups dans la grille"
$ws.Range("G65").Value = "J'ai fini la fonction qui affiche les coups sur la grille. Ça ne marche pas"

# ---------------------------------------------------------------------------
# 7. Row 66 - started debugging with Romain, still in progress (no end time).
# ---------------------------------------------------------------------------
$ws.Range("A66").Value = 44623
$ws.Range("B66").Value = 0.44444444444444442
$ws.Range("D66").Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$ws.Range("E66").Value = "CPNV"
$ws.Range("F66").Value = "Affichage des coups dans la grille"
$ws.Range("G66").Value = "J'ai comencé à débugger ma fonction qui affiche les coups sur la grille"
$ws.Range("H66").Value = "Romain"

# ---------------------------------------------------------------------------
# 8. H55 note grew a link to the ASCII table google search M. Viret showed.
# ---------------------------------------------------------------------------
$ws.Range("H55").Value = "M. Viret https://www.google.com/search?q=ascii+table&rlz=1C1PNJJ_frCH968CH968&source=lnms&tbm=isch&sa=X&ved=2ahUKEwjB0t_xsqn2AhURP-wKHZtWAAIQ_AUoAXoECAEQAw&biw=1744&bih=917&dpr=1.1#imgrc=P0o7oyvJKm4ETM"

# ---------------------------------------------------------------------------
# 9. H49 got re-confirmed as a live hyperlink (same URL as its text).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("H49"), "https://www.thesprucecrafts.com/the-basic-rules-of-battleship-411069") | Out-Null
$ws.Range("H32").Copy() | Out-Null
$ws.Range("H49").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 10. Keep the "big duration" conditional format covering the data rows.
# ---------------------------------------------------------------------------
$durFc = $ws.Cells.FormatConditions.Item(4)
$durFc.ModifyAppliesToRange($ws.Range("D1:D196")) | Out-Null

# ---------------------------------------------------------------------------
# 11. Restore the selection/scroll position left behind by the editing pass.
# ---------------------------------------------------------------------------
$ws.Range("H40").Select() | Out-Null
